$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# Restrict the activity bound to power (EN* -> EN_Z*) since EN* conflicts with
# some UCs for hydro in some regions, and remove the small activity bound
# (AllRegions / I7) for the 2025 row since it caused infeasibilities.
$ws.Range("I7").Value = 0
$ws.Range("L7").Value = "EN_Z*"
$ws.Range("L8").Value = "EN_Z*"
$ws.Range("L9").Value = "EN_Z*"

# Leave a note on the Year cell explaining why it was changed from 2020 to 2025.
$excel.UserName = "Mahmoud Mobir"
$comment = $ws.Range("F7").AddComment("Mahmoud Mobir:" + [char]10 + "12-8-2021" + [char]10 + "This was 2020 but it caused infeasiblities. Made it 2025. ")
$comment.Visible = $false

# Restore selection/view state similar to the saved workbook.
$ws.Range("P11").Select() | Out-Null

$wb.Save() | Out-Null
